$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the style used by the "CLT" rows (G6:G13) off row 6, to clone for the new rows.
# Row 14 mirrors the PJ/Estagio-like block style (same as rows 6-9, fillId 8 -> style index 14 / s=14 for col E, s=15 for G, s=16 for H/I/D)
# Row 15 mirrors CLT block style (rows 10-13) but with a NEW fill (fillId 9) on column G only.

# --- Row 14: José Carlos ---
$ws.Range("A14").Value = "José"
$ws.Range("B14").Value = "Carlos"
$ws.Range("C14").Value = "145.609.357-69"
$ws.Range("E14").Value = "PJ"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = "Desenvolvedor Mobile"

# --- Row 15: Ana Silva ---
$ws.Range("A15").Value = "Ana"
$ws.Range("B15").Value = "Silva"
$ws.Range("C15").Value = "500.000.145-33"
$ws.Range("E15").Value = "CLT"
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = "Desenvolvedor Backend"

# Copy styles from an existing similar row (row 6, the PJ-like block) to row 14
$ws.Range("A6:I6").Copy()
$ws.Range("A14:I14").PasteSpecial(-4122) # xlPasteFormats

# Copy styles from an existing CLT row (row 10) to row 15
$ws.Range("A10:I10").Copy()
$ws.Range("A15:I15").PasteSpecial(-4122) # xlPasteFormats

# Clear the values that shouldn't carry over from the copied formatting (D, H, I are blank in new rows)
$ws.Range("D14").Value = $null
$ws.Range("H14").Value = $null
$ws.Range("I14").Value = $null

$ws.Range("D15").Value = $null
$ws.Range("H15").Value = $null
$ws.Range("I15").Value = $null

# Apply the new unique fill style (fillId 9, centered) to G15 specifically
$ws.Range("G15").Interior.ThemeColor = 7
$ws.Range("G15").Interior.TintAndShade = 0.59999389629810485
$ws.Range("G15").HorizontalAlignment = -4108 # xlCenter

# Update selection / active cell to match final state
$ws.Range("G15").Select()
